$wb = $excel.ActiveWorkbook

# --- "Back-End" sheet (sheet1): row 11 gets a Contributor and its Status
#     moves from "Pending" to "In Progress" (styled like a gray/white
#     "Check Cell"), matching the pattern already used in rows 4, 5 and 10.
$backEnd = $wb.Worksheets.Item(1)

$status = $backEnd.Range("B11")
$status.Value = "In Progress"
$status.Style = "Check Cell"
$status.HorizontalAlignment = -4108
$status.VerticalAlignment = -4108

$contributor = $backEnd.Range("C11")
$contributor.Value = "Rex"
$contributor.Style = "Normal"
$contributor.HorizontalAlignment = -4108
$contributor.VerticalAlignment = -4108

# Move the active selection on "Back-End" like in the saved workbook.
$backEnd.Range("A21").Select()

# --- "UI-UX" sheet (sheet2): only the saved selection moved.
$uiUx = $wb.Worksheets.Item(2)
$uiUx.Range("A10").Select()
